# Add the "purchase_page" worksheet by cloning the structurally-identical
# "redeem_voucher" sheet (same layout/styles), positioning it right before
# "redeem_voucher" (i.e. right after "free_bonus"), then overwriting the
# values with the purchase_page test-case summary numbers.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("redeem_voucher")
$src.Copy($src)

# The freshly inserted copy lands immediately before $src ("redeem_voucher"),
# i.e. it is now the sheet right after "free_bonus".
$newSheet = $wb.Worksheets.Item("redeem_voucher (2)")
$newSheet.Name = "purchase_page"

# Execute / Pass / Fail / Not Tested counters.
$newSheet.Range("E3").Value = 27
$newSheet.Range("E4").Value = 16
$newSheet.Range("E5").Value = 6
$newSheet.Range("E6").Value = 5

# Ref. ID of Failed Test Case.
$newSheet.Range("E8").Value = "SYM-PP-04,`nSYM-PP-08,`nSYM-PP-019,`nSYM-PP-022,`nSYM-PP-024,`nSYM-PP-025"

# Match the taller row height used for the wrapped failure-id text.
$newSheet.Rows.Item(8).RowHeight = 98.25

# Header / title cell (set last so it lands after the failure-id string in
# the shared-string table, matching the authored order).
$newSheet.Range("D2").Value = " Test Case Summary (15-02-24)"

# Keep the active/selected cell on the new sheet the way it was captured.
$newSheet.Range("D2:E2").Select()
